$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38
$ws.Range("H7").Value = 2802
$ws.Range("I7").Value = 2802
$ws.Range("K7").Value = 2802
$ws.Range("M7").Value = -2690
$ws.Range("H14").Value = 2802
$ws.Range("I14").Value = 2802
$ws.Range("K14").Value = 2802
$ws.Range("M14").Value = -2611
$ws.Range("H20").Value = 1270.1666
$ws.Range("I20").Value = 1270.1666
$ws.Range("K20").Value = 1270.1666
$ws.Range("M20").Value = -1040.1666
$ws.Range("H35").Value = 1270.1666
$ws.Range("I35").Value = 1270.1666
$ws.Range("K35").Value = 1270.1666
$ws.Range("M35").Value = -891.1666
$ws.Range("H39").Value = 6442.857
$ws.Range("I39").Value = 2114.8572
$ws.Range("J39").Value = 10770.857
$ws.Range("K39").Value = 6344.571599999999
$ws.Range("L39").Value = 32312.571
$ws.Range("M39").Value = -6048.571599999999
$ws.Range("N39").Value = -32904.571
$ws.Range("H74").Value = 8000
$ws.Range("I74").Value = 8000
$ws.Range("K74").Value = 8000
$ws.Range("M74").Value = -7064
$ws.Range("H77").Value = 8000
$ws.Range("I77").Value = 8000
$ws.Range("K77").Value = 40000
$ws.Range("M77").Value = -35320
$ws.Range("H99").Value = 83335280
$ws.Range("I99").Value = 333333340
$ws.Range("J99").Value = 2596.6667
$ws.Range("K99").Value = 1000000020
$ws.Range("L99").Value = 7790.000100000001
$ws.Range("M99").Value = -999998522
$ws.Range("N99").Value = -10786.0001
$ws.Range("H103").Value = 999.6
$ws.Range("I103").Value = 1024.5
$ws.Range("K103").Value = 3073.5
$ws.Range("M103").Value = -2487.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H61").Value = 2746.3333
$ws.Range("I61").Value = 2746.3333
$ws.Range("K61").Value = 2746.3333
$ws.Range("M61").Value = -2534.3333
$ws.Range("H136").Value = 2746.3333
$ws.Range("I136").Value = 2746.3333
$ws.Range("K136").Value = 8238.999899999999
$ws.Range("M136").Value = -5688.999899999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 10.5
$ws.Range("I8").Value = 10.5
$ws.Range("K8").Value = 10.5
$ws.Range("M8").Value = 129.5
$ws.Range("H99").Value = 4163.769
$ws.Range("I99").Value = 4010.75
$ws.Range("K99").Value = 4010.75
$ws.Range("M99").Value = -2512.75
$ws.Range("H100").Value = 46430
$ws.Range("J100").Value = 46430
$ws.Range("L100").Value = 46430
$ws.Range("N100").Value = -48594
$ws.Range("H134").Value = 4900
$ws.Range("I134").Value = 4900
$ws.Range("K134").Value = 14700
$ws.Range("M134").Value = -12165
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5857.6
$ws.Range("I31").Value = 3929
$ws.Range("K31").Value = 3929
$ws.Range("M31").Value = -3634
$ws.Range("H32").Value = 1847
$ws.Range("I32").Value = 1808.75
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 1808.75
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -1492.75
$ws.Range("N32").Value = -2632
$ws.Range("H34").Value = 5857.6
$ws.Range("I34").Value = 3929
$ws.Range("K34").Value = 3929
$ws.Range("M34").Value = -3727
$ws.Range("H54").Value = 25750
$ws.Range("J54").Value = 31500
$ws.Range("L54").Value = 31500
$ws.Range("N54").Value = -32816
$ws.Range("H86").Value = 111117930
$ws.Range("I86").Value = 142860900
$ws.Range("K86").Value = 142860900
$ws.Range("M86").Value = -142859777
$ws.Range("H89").Value = 111117930
$ws.Range("I89").Value = 142860900
$ws.Range("K89").Value = 714304500
$ws.Range("M89").Value = -714298884

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3904.9844
$ws.Range("J55").Value = 3904.9844
$ws.Range("L55").Value = 11714.9532
$ws.Range("N55").Value = -12068.9532
$ws.Range("H80").Value = 5600
$ws.Range("I80").Value = 5600
$ws.Range("K80").Value = 16800
$ws.Range("M80").Value = -15864
$ws.Range("H83").Value = 5600
$ws.Range("I83").Value = 5600
$ws.Range("K83").Value = 50400
$ws.Range("M83").Value = -45720
$ws.Range("H109").Value = 1913.5
$ws.Range("I109").Value = 1913.5
$ws.Range("K109").Value = 5740.5
$ws.Range("M109").Value = -4700.5
$ws.Range("H131").Value = 3229.9
$ws.Range("I131").Value = 756.25
$ws.Range("J131").Value = 4879
$ws.Range("K131").Value = 2268.75
$ws.Range("L131").Value = 14637
$ws.Range("M131").Value = 2771.25
$ws.Range("N131").Value = -24717

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 410615.38
$ws.Range("I11").Value = 658500
$ws.Range("K11").Value = 658500
$ws.Range("M11").Value = -658361
$ws.Range("H14").Value = 100000
$ws.Range("I14").Value = 100000
$ws.Range("K14").Value = 100000
$ws.Range("M14").Value = -99832
$ws.Range("H132").Value = 1676.1428
$ws.Range("I132").Value = 1146.6
$ws.Range("K132").Value = 3439.8
$ws.Range("M132").Value = -909.7999999999997

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1099.25
$ws.Range("I4").Value = 1099.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1099.25
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -986.25
$ws.Range("N4").ClearContents()
$ws.Range("H28").Value = 1099.25
$ws.Range("I28").Value = 1099.25
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1099.25
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -867.25
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 1099.25
$ws.Range("I37").Value = 1099.25
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1099.25
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -992.25
$ws.Range("N37").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 999
$ws.Range("I9").Value = 999
$ws.Range("K9").Value = 999
$ws.Range("M9").Value = -859
